$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number (e.g. "1.006") must be forced to
# Text format first, otherwise Excel auto-converts the inline string into a
# numeric value and silently drops meaningful trailing/leading zeros
# (e.g. "0.5230" -> 0.523). Flag them up front so the value assignment below
# keeps them as text, matching the original inlineStr cell content.
$textCells = @(
    'D4',
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D12',
    'D13',
    'D15',
    'D16',
    'D18',
    'D19',
    'D20',
    'D21',
    'D24',
    'D25',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D48',
    'D49',
    'D50',
    'D51',
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin price / volume figures scraped by the latest GitHub
# Actions run (cell address -> new value).
$updates = @{
    'D2' = '30.362.71'
    'E2' = '  +2.45%  '
    'D3' = '2.109.82'
    'E3' = '  +1.10%  '
    'D4' = '1.006'
    'E4' = '  -0.38%  '
    'D5' = '345.04'
    'E5' = '  +0.52%  '
    'D6' = '1.005'
    'E6' = '  -0.31%  '
    'D7' = '0.5230'
    'E7' = '  +1.51%  '
    'D8' = '0.4449'
    'E8' = '  +1.39%  '
    'D9' = '54.08'
    'E9' = '  +3.90%  '
    'D10' = '0.09390'
    'E10' = '  +1.79%  '
    'D11' = '1.176'
    'E11' = '  +0.11%  '
    'D12' = '25.19'
    'E12' = '  +0.26%  '
    'B13' = 'Chainlink'
    'C13' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D13' = '8.704'
    'E13' = '  +6.66%  '
    'B14' = 'WrappedEther'
    'C14' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D14' = '2.139.29'
    'E14' = '  +2.11%  '
    'D15' = '6.930'
    'E15' = '  +2.85%  '
    'D16' = '101.96'
    'E16' = '  +2.05%  '
    'E17' = '  +0.81%  '
    'D18' = '1.007'
    'E18' = '  -0.20%  '
    'D19' = '21.36'
    'E19' = '  +1.78%  '
    'D20' = '0.06723'
    'E20' = '  +1.45%  '
    'D21' = '6.309'
    'E21' = '  +2.15%  '
    'E22' = '  -0.13%  '
    'D23' = '30.417.52'
    'E23' = '  +2.42%  '
    'D24' = '12.67'
    'E24' = '  +0.21%  '
    'D25' = '2.320'
    'E25' = '  +0.40%  '
    'D26' = '2.386.31'
    'E26' = '  +2.02%  '
    'D27' = '22.04'
    'E27' = '  +0.95%  '
    'D28' = '2.544'
    'E28' = '  +1.24%  '
    'D29' = '162.28'
    'E29' = '  -0.48%  '
    'D30' = '133.62'
    'E30' = '  +0.97%  '
    'D31' = '1.151'
    'E31' = '  +1.16%  '
    'D32' = '1.761'
    'E32' = '  +8.29%  '
    'D33' = '0.1057'
    'E33' = '  +0.63%  '
    'D34' = '6.859'
    'E34' = '  +13.73%  '
    'D35' = '6.270'
    'E35' = '  +1.57%  '
    'D36' = '3.943'
    'E36' = '  -0.38%  '
    'D37' = '10.61'
    'E37' = '  +2.70%  '
    'D38' = '0.02643'
    'E38' = '  +2.86%  '
    'D39' = '0.06844'
    'E39' = '  +2.06%  '
    'D40' = '0.7086'
    'E40' = '  +3.90%  '
    'D41' = '12.60'
    'E41' = '  +1.22%  '
    'D42' = '1.338'
    'E42' = '  +3.66%  '
    'D43' = '0.2237'
    'E43' = '  +0.00%  '
    'D44' = '0.6867'
    'E44' = '  +3.71%  '
    'D45' = '14.48'
    'E45' = '  +2.38%  '
    'D46' = '2.375'
    'E46' = '  +3.00%  '
    'D48' = '1.394'
    'E48' = '  +19.38%  '
    'D49' = '3.641'
    'E49' = '  +0.81%  '
    'B50' = 'EOS'
    'C50' = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
    'D50' = '1.225'
    'E50' = '  +0.79%  '
    'B51' = 'ThetaToken'
    'C51' = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
    'D51' = '1.202'
    'E51' = '  +8.38%  '
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
